$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 109, shifting rows 109:149 down to 110:150
$ws.Rows("109:109").Insert()

# Populate the newly inserted row 109 with the new Perejil record
$ws.Range("A109").Value = 8
$ws.Range("B109").Value = "Terminal La Palmera de La Serena"
$ws.Range("C109").Value = "Coquimbo"
$ws.Range("D109").Value = 44726
$ws.Range("D109").NumberFormat = $ws.Range("D110").NumberFormat
$ws.Range("E109").Value = 4
$ws.Range("F109").Value = 100112044
$ws.Range("G109").Value = "Perejil"
$ws.Range("H109").Value = "Sin especificar"
$ws.Range("I109").Value = "Primera"
$ws.Range("J109").Value = 3200
$ws.Range("K109").Value = 1500
$ws.Range("L109").Value = 2000
$ws.Range("M109").Value = 1750
$ws.Range("N109").Value = "`$/atado 1 a 1,5 kilos"
$ws.Range("O109").Value = "Provincia del Elquí"
$ws.Range("P109").Value = 1167
$ws.Range("Q109").Value = 1.5
$ws.Range("R109").Value = "Hortaliza"
